$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

function Set-TextValue($addr, $val) {
    # The source cells hold numeric-looking strings as TEXT (not numbers).
    # Temporarily force a text number-format so Excel doesn't silently
    # convert the assigned string into a real number, then drop the
    # number-format override again so the cell's formatting is left as it
    # was originally (only its content/type changes).
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.ClearFormats()
}

# Price (column D) updates
Set-TextValue "D3" "23.13"
Set-TextValue "D4" "5.413"
Set-TextValue "D5" "0.06048"
Set-TextValue "D6" "3.397"
Set-TextValue "D7" "0.8074"
Set-TextValue "D8" "0.9354"
Set-TextValue "D10" "0.07438"
Set-TextValue "D12" "0.03073"
Set-TextValue "D13" "0.09362"
Set-TextValue "D14" "3.943"
Set-TextValue "D15" "0.001595"
Set-TextValue "D16" "0.04844"
Set-TextValue "D18" "0.005386"
Set-TextValue "D19" "0.004165"
Set-TextValue "D20" "0.0009839"
Set-TextValue "D22" "3.650"
Set-TextValue "D23" "6.441"
Set-TextValue "D24" "2.186"
Set-TextValue "D40" "0.03979"
Set-TextValue "D41" "0.006413"
Set-TextValue "D44" "0.006129"
Set-TextValue "D45" "0.00005211"
Set-TextValue "D49" "0.002180"

# Volume(1h) (column E) text updates
$ws.Range("E41").Value = "40KickTokenKICK"
$ws.Range("E48").Value = "47CoinbaseStockTokenCOINBestin24h"
